$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per refreshed crypto data.
# D-column values that look purely numeric get an explicit Text number
# format first so Excel stores them verbatim (matches source "inlineStr" cells)
# instead of silently parsing them into floating-point numbers.
$ws.Range("D2").Value = "25.949.83"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.641.84"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.39"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5082"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06381"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.47"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07782"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "1.637.31"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5461"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "0.0₅7844"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.45"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "26.011.83"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.03"
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.441"
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.962"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.055"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.877"
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.95"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1147"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.883"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.240"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05028"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.260"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.542"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.371"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8992"
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.598"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("D37").Value = "1.132.10"
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5498"
$ws.Range("D39").Value = "0.0₈135"
$ws.Range("E39").Value = "  +15.04%  "
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.618"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8195"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.23"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("D46").Value = "1.779.63"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4533"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.94"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05072"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.007"
$ws.Range("E51").Value = "  +0.40%  "

# Rows 41/42: coin ordering swapped (PaxDollar now ranks above mCoin),
# each with freshly refreshed price/volume figures.
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.005"
$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("B42").Value = "mCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.551"
$ws.Range("E42").Value = "  -0.84%  "
